$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the sheet: "주문등록-양식파일" -> "주문등록"
# ---------------------------------------------------------------------------
$ws.Name = "주문등록"

# ---------------------------------------------------------------------------
# 2. Row 1 = column headers (reordered + corrected wording), Row 2 = sample /
#    description values. The old template had an extra stray column (M) that
#    is no longer used, so it gets dropped later once the new 12-column
#    (A:L) layout is in place.
# ---------------------------------------------------------------------------
$headers = @(
    "자체주문번호",
    "주문자명",
    "주문자전화번호",
    "수령자명",
    "수령자휴대폰번호",
    "우편번호",
    "수령자주소",
    "배송메시지",
    "상품코드",
    "상품명",
    "수량",
    "옵션"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$descriptions = @(
    "쇼핑몰주문번호 또는 자체관리번호",
    "주문자 이름",
    "주문자 연락처",
    "받는분 이름",
    "받는분 휴대폰",
    "우편번호(선택)",
    "전체주소",
    "배송시 요청사항",
    "상품코드",
    "상품명",
    "주문수량",
    "옵션(선택)"
)
for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $descriptions[$i]
}

# ---------------------------------------------------------------------------
# 3. Formatting. Build each finished look (font + fill + border + alignment)
#    once on an out-of-the-way scratch cell, then stamp it onto the real
#    ranges with a single Paste-Special(Formats) so every cell in the row
#    ends up sharing one consolidated style.
# ---------------------------------------------------------------------------
$headerScratch = $ws.Range("N1")
$headerScratch.Interior.Color = 12874308   # RGB(68,114,196) - blue
$headerScratch.Borders.LineStyle = 1       # thin border, all sides
$headerScratch.HorizontalAlignment = -4108 # xlCenter
$headerScratch.VerticalAlignment = -4108   # xlCenter
$headerScratch.Font.Color = 16777215       # RGB(255,255,255) - white
$headerScratch.Font.Bold = $true

$descScratch = $ws.Range("N2")
$descScratch.Interior.Color = 65535        # RGB(255,255,0) - yellow
$descScratch.Borders.LineStyle = 1         # thin border, all sides
$descScratch.HorizontalAlignment = -4108   # xlCenter
$descScratch.VerticalAlignment = -4108     # xlCenter
$descScratch.Font.Color = 6710886          # RGB(102,102,102) - gray
$descScratch.Font.Size = 10

$headerScratch.Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)     # xlPasteFormats

$descScratch.Copy()
$ws.Range("A2:L2").PasteSpecial(-4122)     # xlPasteFormats

$excel.CutCopyMode = $false
$headerScratch.Clear()
$descScratch.Clear()

# ---------------------------------------------------------------------------
# 4. Column widths (characters) and row heights (points)
# ---------------------------------------------------------------------------
$widths = @(20, 12, 15, 12, 15, 10, 40, 25, 15, 25, 8, 15)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 0.83203125
}

$ws.Rows.Item(1).RowHeight = 25
$ws.Rows.Item(2).RowHeight = 20

# ---------------------------------------------------------------------------
# 5. Drop the old trailing column (formerly M, now pushed further right by
#    the scratch work) so the sheet ends cleanly at column L.
# ---------------------------------------------------------------------------
$ws.Columns.Item(13).Delete()

# ---------------------------------------------------------------------------
# 6. Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1    # xlPortrait
$ws.PageSetup.Zoom = 100
